# Daily TGP (terminal gate pricing) refresh: shift the two trading dates forward by
# one day (each old "today" row becomes "yesterday", a new "today" row is added with
# freshly updated Diesel/ULP/PULP/e10 prices) across every state/terminal block on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 46065
$ws.Range("D8").Value = 159.22999999999999
$ws.Range("E8").Value = 149.66999999999999
$ws.Range("F8").Value = 159.66999999999999
$ws.Range("G8").Value = 149.56

# Row 9
$ws.Range("A9").Value = 46065
$ws.Range("D9").Value = 159.22999999999999
$ws.Range("E9").Value = 149.66999999999999
$ws.Range("F9").Value = 159.66999999999999
$ws.Range("G9").Value = 149.56

# Row 10
$ws.Range("A10").Value = 46065
$ws.Range("D10").Value = 160.81
$ws.Range("E10").Value = 152.22999999999999
$ws.Range("F10").Value = 162.22999999999999
$ws.Range("G10").Value = 152.47999999999999

# Row 11
$ws.Range("A11").Value = 46064
$ws.Range("D11").Value = 159.09
$ws.Range("E11").Value = 148.81
$ws.Range("F11").Value = 158.81
$ws.Range("G11").Value = 148.69999999999999

# Row 12
$ws.Range("A12").Value = 46064
$ws.Range("D12").Value = 159.09
$ws.Range("E12").Value = 148.81
$ws.Range("F12").Value = 158.81
$ws.Range("G12").Value = 148.69999999999999

# Row 13
$ws.Range("A13").Value = 46064
$ws.Range("D13").Value = 160.62
$ws.Range("E13").Value = 151.32
$ws.Range("F13").Value = 161.32
$ws.Range("G13").Value = 151.56

# Row 17
$ws.Range("A17").Value = 46065
$ws.Range("D17").Value = 165.07
$ws.Range("E17").Value = 155.68
$ws.Range("F17").Value = 165.68

# Row 18
$ws.Range("A18").Value = 46064
$ws.Range("D18").Value = 164.85
$ws.Range("E18").Value = 154.74
$ws.Range("F18").Value = 164.74

# Row 22
$ws.Range("A22").Value = 46065
$ws.Range("D22").Value = 160.4
$ws.Range("E22").Value = 151.79
$ws.Range("F22").Value = 161.38999999999999
$ws.Range("G22").Value = 153.54

# Row 23
$ws.Range("A23").Value = 46065
$ws.Range("D23").Value = 165.79
$ws.Range("E23").Value = 157.91999999999999
$ws.Range("F23").Value = 167.92

# Row 24
$ws.Range("A24").Value = 46065
$ws.Range("D24").Value = 165.97
$ws.Range("E24").Value = 158.51
$ws.Range("F24").Value = 168.51

# Row 25
$ws.Range("A25").Value = 46065
$ws.Range("D25").Value = 165.97
$ws.Range("E25").Value = 158.03
$ws.Range("F25").Value = 168.03
$ws.Range("G25").Value = 158.88999999999999

# Row 26
$ws.Range("A26").Value = 46065
$ws.Range("D26").Value = 165.58
$ws.Range("E26").Value = 159.63
$ws.Range("F26").Value = 169.63

# Row 27
$ws.Range("A27").Value = 46064
$ws.Range("D27").Value = 160.27000000000001
$ws.Range("E27").Value = 150.93
$ws.Range("F27").Value = 160.53
$ws.Range("G27").Value = 152.68

# Row 28
$ws.Range("A28").Value = 46064
$ws.Range("D28").Value = 165.6
$ws.Range("E28").Value = 157.01
$ws.Range("F28").Value = 167.01

# Row 29
$ws.Range("A29").Value = 46064
$ws.Range("D29").Value = 165.77
$ws.Range("E29").Value = 157.6
$ws.Range("F29").Value = 167.6

# Row 30
$ws.Range("A30").Value = 46064
$ws.Range("D30").Value = 165.77
$ws.Range("E30").Value = 157.12
$ws.Range("F30").Value = 167.12
$ws.Range("G30").Value = 157.97999999999999

# Row 31
$ws.Range("A31").Value = 46064
$ws.Range("D31").Value = 165.38
$ws.Range("E31").Value = 158.72
$ws.Range("F31").Value = 168.72

# Row 35
$ws.Range("A35").Value = 46065
$ws.Range("D35").Value = 159.30000000000001
$ws.Range("E35").Value = 149.91999999999999
$ws.Range("F35").Value = 158.91999999999999

# Row 36
$ws.Range("A36").Value = 46064
$ws.Range("D36").Value = 159.11000000000001
$ws.Range("E36").Value = 149.01
$ws.Range("F36").Value = 158.01

# Row 40
$ws.Range("A40").Value = 46065
$ws.Range("D40").Value = 165.59
$ws.Range("E40").Value = 157.58000000000001
$ws.Range("F40").Value = 167.58

# Row 41
$ws.Range("A41").Value = 46065
$ws.Range("D41").Value = 165.31
$ws.Range("E41").Value = 158
$ws.Range("F41").Value = 168

# Row 42
$ws.Range("A42").Value = 46064
$ws.Range("D42").Value = 165.41
$ws.Range("E42").Value = 156.69
$ws.Range("F42").Value = 166.69

# Row 43
$ws.Range("A43").Value = 46064
$ws.Range("D43").Value = 165.13
$ws.Range("E43").Value = 157.12
$ws.Range("F43").Value = 167.12

# Row 47
$ws.Range("A47").Value = 46065
$ws.Range("D47").Value = 159.59
$ws.Range("E47").Value = 150.62
$ws.Range("F47").Value = 160.62

# Row 48
$ws.Range("A48").Value = 46065
$ws.Range("D48").Value = 159.26
$ws.Range("E48").Value = 150.58000000000001
$ws.Range("F48").Value = 160.58000000000001

# Row 49
$ws.Range("A49").Value = 46064
$ws.Range("D49").Value = 160
$ws.Range("E49").Value = 150.37
$ws.Range("F49").Value = 160.37

# Row 50
$ws.Range("A50").Value = 46064
$ws.Range("D50").Value = 159.66
$ws.Range("E50").Value = 150.32
$ws.Range("F50").Value = 160.32

# Row 54
$ws.Range("A54").Value = 46065
$ws.Range("D54").Value = 174.71
$ws.Range("E54").Value = 165.23
$ws.Range("F54").Value = 175.23

# Row 55
$ws.Range("A55").Value = 46065
$ws.Range("D55").Value = 164.08
$ws.Range("E55").Value = 163.38
$ws.Range("F55").Value = 173.38

# Row 56
$ws.Range("A56").Value = 46065
$ws.Range("D56").Value = 163.95

# Row 57
$ws.Range("A57").Value = 46065
$ws.Range("D57").Value = 164.68
$ws.Range("E57").Value = 157.80000000000001

# Row 58
$ws.Range("A58").Value = 46065
$ws.Range("D58").Value = 160.44999999999999
$ws.Range("E58").Value = 153.69999999999999
$ws.Range("F58").Value = 163.69999999999999

# Row 59
$ws.Range("A59").Value = 46065
$ws.Range("D59").Value = 167.47
$ws.Range("E59").Value = 163.61000000000001

# Row 60
$ws.Range("A60").Value = 46064
$ws.Range("D60").Value = 174.51
$ws.Range("E60").Value = 164.4
$ws.Range("F60").Value = 174.4

# Row 61
$ws.Range("A61").Value = 46064
$ws.Range("D61").Value = 163.89
$ws.Range("E61").Value = 162.44999999999999
$ws.Range("F61").Value = 172.45

# Row 62
$ws.Range("A62").Value = 46064
$ws.Range("D62").Value = 163.76

# Row 63
$ws.Range("A63").Value = 46064
$ws.Range("D63").Value = 164.46
$ws.Range("E63").Value = 156.87

# Row 64
$ws.Range("A64").Value = 46064
$ws.Range("D64").Value = 160.22999999999999
$ws.Range("E64").Value = 152.77000000000001
$ws.Range("F64").Value = 162.77000000000001

# Row 65
$ws.Range("A65").Value = 46064
$ws.Range("D65").Value = 167.25
$ws.Range("E65").Value = 162.75
